$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (admin_level_2 is D, collector_name was E)
# so the new "village" field sits right after admin_level_2, matching the
# shifted columns E..K -> F..L seen in the diff.
$ws.Columns("E").Insert()

# New header cell for the "village" field.
$ws.Range("E1").Value = "village"

# Give column D (admin_level_2) and the new column E (village) the same
# (merged) width, matching the post-edit layout where both share one
# <col> width entry.
$ws.Range("D1:E1").ColumnWidth = 22.666666666666668
